$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header cells: update year labels in D1/E1
$ws.Cells.Item(1, 4).Value = "december 31. 2020"
$ws.Cells.Item(1, 5).Value = "december 31. 2020:1"

# Rows 2-29: update line_num (B), variable label (C), and the two value columns (D,E)
$ws.Cells.Item(2, 2).Value = 3
$ws.Cells.Item(2, 3).Value = " cash and cash equivalents"
$ws.Cells.Item(2, 4).Value = 384344
$ws.Cells.Item(2, 5).Value = 27872

$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = " and respectively"
$ws.Cells.Item(3, 4).Value = 176617
$ws.Cells.Item(3, 5).Value = 148855

$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = " prepaid and other"
$ws.Cells.Item(4, 4).Value = 63224
$ws.Cells.Item(4, 5).Value = 52161

$ws.Cells.Item(5, 2).Value = 8
$ws.Cells.Item(5, 3).Value = " total current assets"
$ws.Cells.Item(5, 4).Value = 624185
$ws.Cells.Item(5, 5).Value = 228888

$ws.Cells.Item(6, 2).Value = 9
$ws.Cells.Item(6, 3).Value = " fixed assets net"
$ws.Cells.Item(6, 4).Value = 628757
$ws.Cells.Item(6, 5).Value = 636153

$ws.Cells.Item(7, 2).Value = 10
$ws.Cells.Item(7, 3).Value = " goodwill"
$ws.Cells.Item(7, 4).Value = 1431967
$ws.Cells.Item(7, 5).Value = 1412873

$ws.Cells.Item(8, 2).Value = 11
$ws.Cells.Item(8, 3).Value = " other intangible asse*s net"
$ws.Cells.Item(8, 4).Value = 274620
$ws.Cells.Item(8, 5).Value = 304673

$ws.Cells.Item(9, 2).Value = 12
$ws.Cells.Item(9, 3).Value = " operating lease right-of-use assets"
$ws.Cells.Item(9, 4).Value = 717821
$ws.Cells.Item(9, 5).Value = 700956

$ws.Cells.Item(10, 2).Value = 13
$ws.Cells.Item(10, 3).Value = " other assets"
$ws.Cells.Item(10, 4).Value = 49298
$ws.Cells.Item(10, 5).Value = 46877

$ws.Cells.Item(11, 2).Value = 14
$ws.Cells.Item(11, 3).Value = " total assets"
$ws.Cells.Item(11, 4).Value = 3726648
$ws.Cells.Item(11, 5).Value = 3330420

$ws.Cells.Item(12, 2).Value = 17
$ws.Cells.Item(12, 3).Value = " current portion of long-term debt"
$ws.Cells.Item(12, 4).Value = 10750
$ws.Cells.Item(12, 5).Value = 10750

$ws.Cells.Item(13, 2).Value = 18
$ws.Cells.Item(13, 3).Value = " accounts payable and accrued expenses"
$ws.Cells.Item(13, 4).Value = 194551
$ws.Cells.Item(13, 5).Value = 167059

$ws.Cells.Item(14, 2).Value = 19
$ws.Cells.Item(14, 3).Value = " current portion of operating lease liabilities"
$ws.Cells.Item(14, 4).Value = 87181
$ws.Cells.Item(14, 5).Value = 83123

$ws.Cells.Item(15, 2).Value = 20
$ws.Cells.Item(15, 3).Value = " deferred"
$ws.Cells.Item(15, 4).Value = 197939
$ws.Cells.Item(15, 5).Value = 191117

$ws.Cells.Item(16, 2).Value = 22
$ws.Cells.Item(16, 3).Value = " other current liabilities"
$ws.Cells.Item(16, 4).Value = 40393
$ws.Cells.Item(16, 5).Value = 31241

$ws.Cells.Item(17, 2).Value = 23
$ws.Cells.Item(17, 3).Value = " total current liabilities"
$ws.Cells.Item(17, 4).Value = 530814
$ws.Cells.Item(17, 5).Value = 483290

$ws.Cells.Item(18, 2).Value = 24
$ws.Cells.Item(18, 3).Value = " long-iciu debt net"
$ws.Cells.Item(18, 4).Value = 1020137
$ws.Cells.Item(18, 5).Value = 1028049

$ws.Cells.Item(19, 2).Value = 25
$ws.Cells.Item(19, 3).Value = " operating lease liabilities"
$ws.Cells.Item(19, 4).Value = 729754
$ws.Cells.Item(19, 5).Value = 685910

$ws.Cells.Item(20, 2).Value = 26
$ws.Cells.Item(20, 3).Value = " other long-term liabilities"
$ws.Cells.Item(20, 4).Value = 105980
$ws.Cells.Item(20, 5).Value = 92865

$ws.Cells.Item(21, 2).Value = 27
$ws.Cells.Item(21, 3).Value = " deferred"
$ws.Cells.Item(21, 4).Value = 10215
$ws.Cells.Item(21, 5).Value = 10098

$ws.Cells.Item(22, 2).Value = 29
$ws.Cells.Item(22, 3).Value = " deferred income taxes"
$ws.Cells.Item(22, 4).Value = 45951
$ws.Cells.Item(22, 5).Value = 58940

$ws.Cells.Item(23, 2).Value = 30
$ws.Cells.Item(23, 3).Value = " total liabilities"
$ws.Cells.Item(23, 4).Value = 2442851
$ws.Cells.Item(23, 5).Value = 2359152

$ws.Cells.Item(24, 2).Value = 36
$ws.Cells.Item(24, 3).Value = " issued and outstanding at december and respectively"
$ws.Cells.Item(24, 4).Value = 60
$ws.Cells.Item(24, 5).Value = 58

$ws.Cells.Item(25, 2).Value = 37
$ws.Cells.Item(25, 3).Value = " additional paid-in capital"
$ws.Cells.Item(25, 4).Value = 910304
$ws.Cells.Item(25, 5).Value = 648031

$ws.Cells.Item(26, 2).Value = 38
$ws.Cells.Item(26, 3).Value = " accumulated other comprehensive loss"
$ws.Cells.Item(26, 4).Value = -27069
$ws.Cells.Item(26, 5).Value = -50331

$ws.Cells.Item(27, 2).Value = 39
$ws.Cells.Item(27, 3).Value = " retained earnings"
$ws.Cells.Item(27, 4).Value = 400502
$ws.Cells.Item(27, 5).Value = 373510

$ws.Cells.Item(28, 2).Value = 40
$ws.Cells.Item(28, 3).Value = " total stockholders’ equity"
$ws.Cells.Item(28, 4).Value = 1283797
$ws.Cells.Item(28, 5).Value = 971268

$ws.Cells.Item(29, 2).Value = 41
$ws.Cells.Item(29, 3).Value = " total liabilities and stockholders’ equity see accompanying notes to consolidated financial statements."
$ws.Cells.Item(29, 4).Value = 3726648
$ws.Cells.Item(29, 5).Value = 3330420

# Rows 30-33: remove the line-item rows entirely (page_num/line_num/variable cleared,
# and the two value columns cleared back to blank placeholders)
$ws.Range("A30:E33").ClearContents()

